$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.367.44"
$ws.Range("E2").Value = "  +4.42%  "

$ws.Range("D3").Value = "4.037.74"
$ws.Range("E3").Value = "  +3.63%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.28"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.97"
$ws.Range("E6").Value = "  +3.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("E7").Value = "  +3.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.744"
$ws.Range("E9").Value = "  +3.66%  "

$ws.Range("E10").Value = "  +2.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000334"
$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.37"
$ws.Range("E12").Value = "  +13.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.91"
$ws.Range("E13").Value = "  +7.05%  "

$ws.Range("D14").Value = "4.685.23"
$ws.Range("E14").Value = "  +3.79%  "

$ws.Range("D15").Value = "4.029.19"
$ws.Range("E15").Value = "  +3.68%  "

$ws.Range("E16").Value = "  +8.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.27"
$ws.Range("E17").Value = "  +3.02%  "

$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.133"
$ws.Range("E19").Value = "  -1.85%  "

$ws.Range("D20").Value = "72.343.52"
$ws.Range("E20").Value = "  +4.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.52"
$ws.Range("E21").Value = "  +3.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "101.29"
$ws.Range("E22").Value = "  +15.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.57"
$ws.Range("E23").Value = "  +6.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.75"
$ws.Range("E24").Value = "  +4.26%  "

$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.87"
$ws.Range("E26").Value = "  +2.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.28"
$ws.Range("E27").Value = "  +7.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.54"
$ws.Range("E28").Value = "  +3.84%  "

$ws.Range("E29").Value = "  +10.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.57"
$ws.Range("E30").Value = "  +3.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "694.06"
$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("E32").Value = "  +3.36%  "

$ws.Range("E33").Value = "  +17.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "68.11"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("B35").Value = "TheGraph"
$ws.Range("C35").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.444"
$ws.Range("E35").Value = "  +1.88%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0889"
$ws.Range("E36").Value = "  +7.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.65"
$ws.Range("E37").Value = "  +4.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.64"
$ws.Range("E38").Value = "  +22.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.153"
$ws.Range("E39").Value = "  +3.31%  "

$ws.Range("E40").Value = "  -0.30%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0491"
$ws.Range("E42").Value = "  +2.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.18"
$ws.Range("E43").Value = "  +6.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("E44").Value = "  +2.11%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.148"
$ws.Range("E45").Value = "  +5.71%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.50"
$ws.Range("E46").Value = "  +5.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.11"
$ws.Range("E47").Value = "  +2.90%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.09"
$ws.Range("E48").Value = "  +9.11%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000275"
$ws.Range("E49").Value = "  +23.48%  "

$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("E51").Value = "  +0.21%  "
